$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Triangular): turn B3:N3 into one shared-formula run (was B3 own + C3:N3 own) ---
$ws.Range("B3:N3").Formula = "=B2*(B2+1)/2"

# --- New row 9: Five ---
$ws.Range("A9").Value2 = "Five"
$ws.Range("B9").Formula = "=B2*5"
$ws.Range("C9").Formula = "=C2*5"
$ws.Range("D9:N9").Formula = "=D2*5"

# --- New row 10: Five Sum ---
$ws.Range("A10").Value2 = "Five Sum"
$ws.Range("B10").Formula = "=B9"
$ws.Range("C10").Formula = "=B10+C9"
$ws.Range("D10").Formula = "=C10+D9"
$ws.Range("E10:N10").Formula = "=D10+E9"

# --- New row 11: Six ---
$ws.Range("A11").Value2 = "Six"
$ws.Range("B11").Formula = "=B2*6"
$ws.Range("C11").Formula = "=C2*6"
$ws.Range("D11:N11").Formula = "=D2*6"

# --- New row 12: Six Sum ---
$ws.Range("A12").Value2 = "Six Sum"
$ws.Range("B12").Formula = "=B11"
$ws.Range("C12").Formula = "=B12+C11"
$ws.Range("D12:N12").Formula = "=C12+D11"

# --- New row 13: Seven ---
$ws.Range("A13").Value2 = "Seven"
$ws.Range("B13").Formula = "=B2*7"
$ws.Range("C13").Formula = "=C2*7"
$ws.Range("D13:N13").Formula = "=D2*7"

# --- New row 14: Seven Sum ---
$ws.Range("A14").Value2 = "Seven Sum"
$ws.Range("B14").Formula = "=B13"
$ws.Range("C14").Formula = "=B14+C13"
$ws.Range("D14:N14").Formula = "=C14+D13"

# --- Selection moves to S30 (matches the new selection recorded in the saved file) ---
$ws.Range("S30").Select()
